$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.033.24'
$ws.Range('D3').Value = '1.668.40'
$ws.Range('E3').Value = '  -1.35%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '216.88'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.45%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5078'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.58%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.004'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('E8').Value = '  -0.39%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06391'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.27%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '21.80'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.17%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07450'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.33%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.667.95'
$ws.Range('E12').Value = '  -1.52%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.513'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.14%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.5825'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.86%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.000008561'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.01%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.35'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.38%  '
$ws.Range('D17').Value = '26.120.22'
$ws.Range('E17').Value = '  -1.80%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '4.936'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.13%  '
$ws.Range('E19').Value = '  -0.12%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '10.79'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.27%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '191.22'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.47%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.193'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.08%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.005'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '144.65'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '7.623'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.67%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1199'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.35%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '15.67'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.70%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.06625'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +15.65%  '
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.317'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.82%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.547'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.73%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.513'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.657'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.88%  '
$ws.Range('E34').Value = '  -0.36%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.6130'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.45%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.370'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.687'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.40%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.299'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +7.95%  '
$ws.Range('D39').Value = '1.095.10'
$ws.Range('E39').Value = '  +0.30%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01599'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.30%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.8712'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.25%  '
$ws.Range('E42').Value = '  +0.32%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '101.13'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.50%  '
$ws.Range('D44').Value = '1.816.79'
$ws.Range('E44').Value = '  -1.60%  '
$ws.Range('E45').Value = '  -2.56%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '56.42'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.16%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.012'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.72%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '8.075'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.12%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.05230'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('E50').Value = '  -0.78%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '6.034'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +4.24%  '
